$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "Eval" -> "EVAL"
$ws.Range("B6").Value = "EVAL"

# Rows 9 & 10: swap the IMPORT row and the "model name" row
# (also updates the "level of imporance" column values accordingly)
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "IMPORT"
$ws.Range("C9").Value = "any specific to imports, libraries, models, data"

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = """model name"""
$ws.Range("C10").Value = "specific to this model"

# Update the selection shown in the sheet view
$ws.Range("A11").Select()
